$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "simov"
$ws.Range("B2").Value = "vc"
$ws.Range("C2").Value = "quer "
$ws.Range("D2").Value = "bbbbbb"
$ws.Range("G2").Value = "naooo"
